$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    # MatchCase=$true so e.g. "Multiple power-up spin options" does not also
    # clobber the lower-cased "multiple power-up spin options" phrase that
    # appears inside the closing summary sentence.
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# Closing italic summary paragraph (replace first - it shares a phrase with
# the bullet list below, so doing it first avoids any case-insensitive
# cross-match issues)
Replace-AllText "Experience the beautiful wilderness of North America in Buffalo Bounty, an online slot game featuring free spins and multiple power-up spin options." "Experience the Wild West and play Buffalo Bounty slot game for free. Enjoy stunning graphics and thrilling features."

# Title (appears twice: Heading1 title + bold closing line)
Replace-AllText "Play Buffalo Bounty Free: Game Review & Features" "Play Buffalo Bounty Free - Exciting Online Slot Game"

# "What we like" bullet list
Replace-AllText "Beautiful graphics of the North American wilderness" "Beautiful graphics depicting scenic landscapes"
Replace-AllText "Free spins feature triggered by Scatter symbol" "Exciting free spins feature with increasing number of spins"
Replace-AllText "Multiple power-up spin options" "Power-up spin options for enhanced gameplay"
Replace-AllText "High-paying buffalo symbol worth up to 15x total bet" "Immersive experience for experienced gamblers"

# "What we don't like" bullet list
Replace-AllText "Only 10 paylines" "Limited number of paylines (10)"
